# Auto-generated edit script
# Applies cached-value updates to Sheets/Behemoth_Profits.xlsx per commit diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1518.6923
$ws.Range("I19").Value = 1072.8334
$ws.Range("K19").Value = 1072.8334
$ws.Range("M19").Value = -897.8334
$ws.Range("H64").Value = 4788
$ws.Range("I64").Value = 3675
$ws.Range("K64").Value = 3675
$ws.Range("M64").Value = -3427
$ws.Range("H67").Value = 4788
$ws.Range("I67").Value = 3675
$ws.Range("K67").Value = 3675
$ws.Range("M67").Value = -2817
$ws.Range("H97").Value = 2517.875
$ws.Range("J97").Value = 2517.875
$ws.Range("L97").Value = 7553.625
$ws.Range("N97").Value = -8545.625
$ws.Range("H105").Value = 45825.168
$ws.Range("I105").Value = 32471
$ws.Range("J105").Value = 48496
$ws.Range("K105").Value = 32471
$ws.Range("L105").Value = 48496
$ws.Range("M105").Value = -28977
$ws.Range("N105").Value = -55484
$ws.Range("H111").Value = 568.6667
$ws.Range("I111").Value = 568.6667
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 1706.0001
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 1360.9999
$ws.Range("N111").ClearContents()
$ws.Range("H135").Value = 715.6087
$ws.Range("I135").Value = 366.73685
$ws.Range("J135").Value = 2372.75
$ws.Range("K135").Value = 3300.63165
$ws.Range("L135").Value = 21354.75
$ws.Range("M135").Value = -765.6316500000003
$ws.Range("N135").Value = -26424.75
$ws.Range("H138").Value = 5366.855
$ws.Range("I138").Value = 4876.8335
$ws.Range("J138").Value = 5419.357
$ws.Range("K138").Value = 14630.5005
$ws.Range("L138").Value = 16258.071
$ws.Range("M138").Value = -9490.500499999998
$ws.Range("N138").Value = -26538.071

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 44171.02
$ws.Range("I32").Value = 38790.324
$ws.Range("K32").Value = 38790.324
$ws.Range("M32").Value = -38503.324
$ws.Range("H61").Value = 8625231
$ws.Range("I61").Value = 3825.8845
$ws.Range("K61").Value = 3825.8845
$ws.Range("M61").Value = -3613.8845
$ws.Range("H136").Value = 8625231
$ws.Range("I136").Value = 3825.8845
$ws.Range("K136").Value = 11477.6535
$ws.Range("M136").Value = -8927.6535
$ws.Range("H139").Value = 86998.8
$ws.Range("J139").Value = 94998.5
$ws.Range("L139").Value = 94998.5
$ws.Range("N139").Value = -105278.5

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 46805.11
$ws.Range("I82").Value = 6376
$ws.Range("J82").Value = 127663.336
$ws.Range("K82").Value = 6376
$ws.Range("L82").Value = 127663.336
$ws.Range("M82").Value = -5993
$ws.Range("N82").Value = -128429.336
$ws.Range("H85").Value = 46805.11
$ws.Range("I85").Value = 6376
$ws.Range("J85").Value = 127663.336
$ws.Range("K85").Value = 6376
$ws.Range("L85").Value = 127663.336
$ws.Range("M85").Value = -5050
$ws.Range("N85").Value = -130315.336
$ws.Range("H134").Value = 4388283
$ws.Range("I134").Value = 2155.647
$ws.Range("K134").Value = 6466.941
$ws.Range("M134").Value = -3931.941

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20035.2
$ws.Range("I31").Value = 26833.166
$ws.Range("J31").Value = 15503.223
$ws.Range("K31").Value = 26833.166
$ws.Range("L31").Value = 15503.223
$ws.Range("M31").Value = -26538.166
$ws.Range("N31").Value = -16093.223
$ws.Range("H34").Value = 20035.2
$ws.Range("I34").Value = 26833.166
$ws.Range("J34").Value = 15503.223
$ws.Range("K34").Value = 26833.166
$ws.Range("L34").Value = 15503.223
$ws.Range("M34").Value = -26631.166
$ws.Range("N34").Value = -15907.223
$ws.Range("H47").Value = 5000
$ws.Range("J47").Value = 5000
$ws.Range("L47").Value = 5000
$ws.Range("N47").Value = -6132
$ws.Range("H62").Value = 3019
$ws.Range("I62").Value = 2603.25
$ws.Range("J62").Value = 3850.5
$ws.Range("K62").Value = 2603.25
$ws.Range("L62").Value = 3850.5
$ws.Range("M62").Value = -1979.25
$ws.Range("N62").Value = -5098.5
$ws.Range("H65").Value = 3019
$ws.Range("I65").Value = 2603.25
$ws.Range("J65").Value = 3850.5
$ws.Range("K65").Value = 13016.25
$ws.Range("L65").Value = 19252.5
$ws.Range("M65").Value = -9896.25
$ws.Range("N65").Value = -25492.5
$ws.Range("H99").Value = 3880.3333
$ws.Range("I99").Value = 3810.6667
$ws.Range("K99").Value = 3810.6667
$ws.Range("M99").Value = -2312.6667
$ws.Range("H126").Value = 3880.3333
$ws.Range("I126").Value = 3810.6667
$ws.Range("K126").Value = 11432.0001
$ws.Range("M126").Value = -8962.000100000001
$ws.Range("H134").Value = 4277.6665
$ws.Range("I134").Value = 2317.2
$ws.Range("K134").Value = 6951.599999999999
$ws.Range("M134").Value = -4416.599999999999

# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1707.2693
$ws.Range("I5").Value = 1071.8235
$ws.Range("K5").Value = 3215.4705
$ws.Range("M5").Value = -3103.4705
$ws.Range("H122").Value = 2509.75
$ws.Range("J122").Value = 2602.3333
$ws.Range("L122").Value = 23420.9997
$ws.Range("N122").Value = -28320.9997
$ws.Range("H135").Value = 1707.2693
$ws.Range("I135").Value = 1071.8235
$ws.Range("K135").Value = 9646.4115
$ws.Range("M135").Value = -7111.4115

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 142751.12
$ws.Range("I19").Value = 141715.72
$ws.Range("K19").Value = 141715.72
$ws.Range("M19").Value = -141427.72
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 60.666668
$ws.Range("I46").Value = 60.666668
$ws.Range("K46").Value = 60.666668
$ws.Range("M46").Value = 95.333332
$ws.Range("H54").Value = 14750
$ws.Range("I54").Value = 10000
$ws.Range("J54").Value = 19500
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = 19500
$ws.Range("M54").Value = -9610
$ws.Range("N54").Value = -20280
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H80").Value = 5206.5
$ws.Range("I80").Value = 4349.4
$ws.Range("J80").Value = 6635
$ws.Range("K80").Value = 4349.4
$ws.Range("L80").Value = 6635
$ws.Range("M80").Value = -3351.4
$ws.Range("N80").Value = -8631
$ws.Range("H83").Value = 5206.5
$ws.Range("I83").Value = 4349.4
$ws.Range("J83").Value = 6635
$ws.Range("K83").Value = 21747
$ws.Range("L83").Value = 33175
$ws.Range("M83").Value = -16755
$ws.Range("N83").Value = -43159
$ws.Range("H126").Value = 4690
$ws.Range("I126").Value = 4133.3335
$ws.Range("J126").Value = 4928.5713
$ws.Range("K126").Value = 12400.0005
$ws.Range("L126").Value = 14785.7139
$ws.Range("M126").Value = -9930.000499999998
$ws.Range("N126").Value = -19725.7139
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12862.333
$ws.Range("I132").Value = 16349.5
$ws.Range("J132").Value = 5888
$ws.Range("K132").Value = 49048.5
$ws.Range("L132").Value = 17664
$ws.Range("M132").Value = -46518.5
$ws.Range("N132").Value = -22724
$ws.Range("H136").Value = 57335.82
$ws.Range("I136").Value = 25861.8
$ws.Range("K136").Value = 77585.39999999999
$ws.Range("M136").Value = -75035.39999999999

# --- Sheet: WVR (index 8) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2541.2856
$ws.Range("I81").Value = 2498.3333
$ws.Range("K81").Value = 4996.6666
$ws.Range("M81").Value = -3935.6666
$ws.Range("H84").Value = 2541.2856
$ws.Range("I84").Value = 2498.3333
$ws.Range("K84").Value = 24983.333
$ws.Range("M84").Value = -19679.333
$ws.Range("H132").Value = 3098827
$ws.Range("I132").Value = 6075.625
$ws.Range("J132").Value = 7597374.5
$ws.Range("K132").Value = 18226.875
$ws.Range("L132").Value = 22792123.5
$ws.Range("M132").Value = -15696.875
$ws.Range("N132").Value = -22797183.5
$ws.Range("H136").Value = 1548963.5
$ws.Range("I136").Value = 4195.727
$ws.Range("J136").Value = 3976455.8
$ws.Range("K136").Value = 12587.181
$ws.Range("L136").Value = 11929367.4
$ws.Range("M136").Value = -10037.181
$ws.Range("N136").Value = -11934467.4
